$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.215.67"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").Value = "2.624.58"
$ws.Range("E3").Value = "  +1.13%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'522.56"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("D6").Value = "'149.03"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.83%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -4.40%  "
$ws.Range("D9").Value = "2.628.81"
$ws.Range("E9").Value = "  +0.92%  "
$ws.Range("E10").Value = "  -5.21%  "
$ws.Range("D11").Value = "'0.106"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.57%  "
$ws.Range("E12").Value = "  -1.44%  "
$ws.Range("E13").Value = "  -0.63%  "
$ws.Range("D14").Value = "3.080.03"
$ws.Range("E14").Value = "  +0.92%  "
$ws.Range("D15").Value = "60.203.76"
$ws.Range("E15").Value = "  -0.48%  "
$ws.Range("D16").Value = "'21.23"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.19%  "
$ws.Range("E17").Value = "  -1.52%  "
$ws.Range("D18").Value = "2.621.57"
$ws.Range("E18").Value = "  +0.81%  "
$ws.Range("E19").Value = "  -2.27%  "
$ws.Range("D20").Value = "'340.94"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.92%  "
$ws.Range("E21").Value = "  -1.23%  "
$ws.Range("E22").Value = "  -1.57%  "
$ws.Range("D23").Value = "'0.994"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("D24").Value = "'60.66"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.50%  "
$ws.Range("E25").Value = "  -1.90%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  -2.18%  "
$ws.Range("D28").Value = "0.0₃0810"
$ws.Range("E28").Value = "  -3.78%  "
$ws.Range("E29").Value = "  -3.75%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.59"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").Value = "'6.00"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -5.03%  "
$ws.Range("E33").Value = "  -2.13%  "
$ws.Range("D34").Value = "'150.40"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.15%  "
$ws.Range("D35").Value = "'3.96"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -4.35%  "
$ws.Range("D36").Value = "'0.919"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.09%  "
$ws.Range("E37").Value = "  -4.79%  "
$ws.Range("D38").Value = "'0.866"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +3.23%  "
$ws.Range("D39").Value = "'36.46"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.32%  "
$ws.Range("E40").Value = "  -2.87%  "
$ws.Range("E41").Value = "  -4.18%  "
$ws.Range("D42").Value = "'289.34"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.36%  "
$ws.Range("D43").Value = "'0.626"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.54%  "
$ws.Range("D44").Value = "'0.100"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.12%  "
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").Value = "'0.0547"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.10%  "
$ws.Range("D47").Value = "'19.51"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("E48").Value = "  +0.91%  "
$ws.Range("D49").Value = "'0.0232"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.87%  "
$ws.Range("D50").Value = "'4.70"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.37%  "
$ws.Range("D51").Value = "1.967.71"
$ws.Range("E51").Value = "  +0.04%  "
